# Update workbook: "Automatic update of files." commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (data rows start at 2)
$lastRow = $ws.UsedRange.Rows.Count

# Column C holds a "Förändrad" (changed) timestamp that was refreshed for every row
# from 45182 to 45184 during this automatic update.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value = 45184
    }
}

# Row 2 specific updates: Signalarter (I2) and Alla arter (Q2) counts incremented,
# and the species list (R2) gained a new entry "Ramaria neoformosa".
$ws.Range("I2").Value = 8
$ws.Range("Q2").Value = 21

$speciesList = "Bombmurkla", "Borsttagging", "Gräddticka", "Knärot", "Garnlav", `
    "Motaggsvamp", "Orange taggsvamp", "Skrovlig taggsvamp", "Spillkråka", `
    "Svartvit taggsvamp", "Tallticka", "Tretåig hackspett", "Utter", "Bronshjon", `
    "Dropptaggsvamp", "Gullgröppa", "Mindre märgborre", "Plattlummer", `
    "Ramaria neoformosa", "Tallfingersvamp", "Vedticka"

$ws.Range("R2").Value = [string]::Join("`r`n", $speciesList)

# Setting the longer, wrapped text above causes Excel to auto-fit the row height;
# restore the original explicit row height so the row's formatting is unaffected.
$ws.Rows.Item(2).RowHeight = 15
